# Applies the "adjust for difference instance" edit:
# - Column B (Flight No.) becomes numeric, matching the row's ID (column A) value
# - Columns C/D (Arrival/Departure Time) get new time-of-day values
# - Column E (Location) gets new location codes
# - Column F (Type) gets new values
# Columns G (Delay Avg) and H (Delay Var) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 2;  B = 1;  C = 0.00625;               D = 0.02777777777777778; E = "E5";  F = 2 },
    @{ Row = 3;  B = 2;  C = 0.03958333333333333;   D = 0.06597222222222222; E = "E26"; F = 2 },
    @{ Row = 4;  B = 3;  C = 0.0375;                D = 0.0625;              E = "D30"; F = 1 },
    @{ Row = 5;  B = 4;  C = 0.0125;                D = 0.03472222222222222; E = "A8";  F = 1 },
    @{ Row = 6;  B = 5;  C = 0.01388888888888889;   D = 0.03888888888888889; E = "D33"; F = 1 },
    @{ Row = 7;  B = 6;  C = 0.03125;               D = 0.05694444444444444; E = "B6";  F = 1 },
    @{ Row = 8;  B = 7;  C = 0.002777777777777778;  D = 0.02430555555555556; E = "A14"; F = 2 },
    @{ Row = 9;  B = 8;  C = 0.02222222222222222;   D = 0.04861111111111111; E = "E28"; F = 1 },
    @{ Row = 10; B = 9;  C = 0.006944444444444444;  D = 0.02986111111111111; E = "D37"; F = 1 },
    @{ Row = 11; B = 10; C = 0.02013888888888889;   D = 0.04166666666666666; E = "A11"; F = 3 },
    @{ Row = 12; B = 11; C = 0.02361111111111111;   D = 0.04930555555555555; E = "A10"; F = 3 },
    @{ Row = 13; B = 12; C = 0.004166666666666667;  D = 0.02847222222222222; E = "F41"; F = 2 },
    @{ Row = 14; B = 13; C = 0.03472222222222222;   D = 0.06180555555555556; E = "A1";  F = 2 },
    @{ Row = 15; B = 14; C = 0.01597222222222222;   D = 0.03958333333333333; E = "F54"; F = 3 },
    @{ Row = 16; B = 15; C = 0.03333333333333333;   D = 0.05625;             E = "C2";  F = 2 },
    @{ Row = 17; B = 16; C = 0.0006944444444444445; D = 0.02708333333333333; E = "D45"; F = 1 },
    @{ Row = 18; B = 17; C = 0.01944444444444444;   D = 0.04097222222222222; E = "C18"; F = 2 },
    @{ Row = 19; B = 18; C = 0.02222222222222222;   D = 0.04583333333333333; E = "F59"; F = 3 },
    @{ Row = 20; B = 19; C = 0.03472222222222222;   D = 0.06041666666666667; E = "D46"; F = 1 },
    @{ Row = 21; B = 20; C = 0.03263888888888889;   D = 0.05763888888888889; E = "E27"; F = 1 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
